$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '65.673.03'
$ws.Range("E2").Value = '  +2.93%  '
$ws.Range("D3").Value = '2.664.71'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '606.47'
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("D6").Value = '158.58'
$ws.Range("E6").Value = '  +4.90%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("D9").Value = '0.124'
$ws.Range("E9").Value = '  +8.24%  '
$ws.Range("D10").Value = '0.406'
$ws.Range("E10").Value = '  +2.55%  '
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("D13").Value = '29.82'
$ws.Range("E13").Value = '  +6.90%  '
$ws.Range("D14").Value = '0.0000195'
$ws.Range("E14").Value = '  +15.41%  '
$ws.Range("D15").Value = '3.149.69'
$ws.Range("E15").Value = '  +1.61%  '
$ws.Range("D16").Value = '65.404.95'
$ws.Range("E16").Value = '  +2.63%  '
$ws.Range("D17").Value = '2.681.11'
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").Value = '12.78'
$ws.Range("E18").Value = '  +4.95%  '
$ws.Range("D19").Value = '4.90'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("D20").Value = '360.53'
$ws.Range("E20").Value = '  +3.75%  '
$ws.Range("D21").Value = '7.38'
$ws.Range("E21").Value = '  +5.21%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '68.98'
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").Value = '1.71'
$ws.Range("E24").Value = '  +1.14%  '
$ws.Range("D25").Value = '9.58'
$ws.Range("E25").Value = '  +4.97%  '
$ws.Range("D26").Value = '0.0000106'
$ws.Range("E26").Value = '  +17.45%  '
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("D28").Value = '8.25'
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("E29").Value = '  +1.79%  '
$ws.Range("D30").Value = '2.21'
$ws.Range("E30").Value = '  +7.02%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '537.84'
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("D33").Value = '1.86'
$ws.Range("E33").Value = '  +2.20%  '
$ws.Range("D34").Value = '5.64'
$ws.Range("E34").Value = '  +5.09%  '
$ws.Range("D35").Value = '6.37'
$ws.Range("E35").Value = '  +4.31%  '
$ws.Range("E36").Value = '  +4.30%  '
$ws.Range("D37").Value = '20.70'
$ws.Range("E37").Value = '  +3.67%  '
$ws.Range("D38").Value = '2.02'
$ws.Range("E38").Value = '  +2.44%  '
$ws.Range("D39").Value = '162.94'
$ws.Range("E39").Value = '  -0.70%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = '42.45'
$ws.Range("E42").Value = '  +6.42%  '
$ws.Range("D43").Value = '167.06'
$ws.Range("E43").Value = '  -0.62%  '
$ws.Range("D44").Value = '4.19'
$ws.Range("E44").Value = '  +2.71%  '
$ws.Range("D45").Value = '2.37'
$ws.Range("E45").Value = '  +8.23%  '
$ws.Range("D46").Value = '0.0614'
$ws.Range("E46").Value = '  +5.15%  '
$ws.Range("D47").Value = '23.20'
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("D48").Value = '0.663'
$ws.Range("E48").Value = '  +4.10%  '
$ws.Range("D49").Value = '0.0266'
$ws.Range("E49").Value = '  +5.58%  '
$ws.Range("E50").Value = '  +2.31%  '
$ws.Range("D51").Value = '19.87'
$ws.Range("E51").Value = '  +3.09%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"

